# Auto-generated edit: appends new sensor-log rows to the PIR and
# Humidity sheets (new readings logged on 2026-01-30, 18:20-18:21).
$wb = $excel.ActiveWorkbook

# Each inner array is: Date, Timestamp, Hour, Location, Value, Status
$pirRows = New-Object System.Collections.ArrayList
[void]$pirRows.Add(@('29', '2026-01-30', '18:20:46', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('30', '2026-01-30', '18:20:46', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('31', '2026-01-30', '18:20:51', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('32', '2026-01-30', '18:20:56', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('33', '2026-01-30', '18:21:01', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('34', '2026-01-30', '18:21:06', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('35', '2026-01-30', '18:21:11', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('36', '2026-01-30', '18:21:16', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('37', '2026-01-30', '18:21:21', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('38', '2026-01-30', '18:21:26', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('39', '2026-01-30', '18:21:31', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('40', '2026-01-30', '18:21:37', '18:00', 'Bathroom', 'No Motion', 'Inactive'))
[void]$pirRows.Add(@('41', '2026-01-30', '18:21:41', '18:00', 'Bathroom', 'No Motion', 'Inactive'))

$humRows = New-Object System.Collections.ArrayList
[void]$humRows.Add(@('23', '2026-01-30', '18:20:46', '18:00', 'Bathroom', '86.6%', 'Active'))
[void]$humRows.Add(@('24', '2026-01-30', '18:20:46', '18:00', 'Bathroom', '86.7%', 'Active'))
[void]$humRows.Add(@('25', '2026-01-30', '18:20:52', '18:00', 'Bathroom', '86.7%', 'Active'))
[void]$humRows.Add(@('26', '2026-01-30', '18:21:02', '18:00', 'Bathroom', '85.2%', 'Active'))
[void]$humRows.Add(@('27', '2026-01-30', '18:21:07', '18:00', 'Bathroom', '85.8%', 'Active'))
[void]$humRows.Add(@('28', '2026-01-30', '18:21:12', '18:00', 'Bathroom', '86.7%', 'Active'))
[void]$humRows.Add(@('29', '2026-01-30', '18:21:22', '18:00', 'Bathroom', '86.7%', 'Active'))
[void]$humRows.Add(@('30', '2026-01-30', '18:21:32', '18:00', 'Bathroom', '86.7%', 'Active'))
[void]$humRows.Add(@('31', '2026-01-30', '18:21:42', '18:00', 'Bathroom', '86.7%', 'Active'))

function Add-LogRows([object]$ws, [object]$rows) {
  foreach ($entry in $rows) {
    $r = [int]$entry[0]
    for ($c = 1; $c -le 6; $c++) {
      $cell = $ws.Cells.Item($r, $c)
      # Force text so date/time-looking strings ("2026-01-30", "18:20:46")
      # are stored verbatim instead of being auto-converted to date/time
      # serial numbers, then drop back to the default "Normal" style so no
      # stray number-format style sticks to the cell.
      $cell.NumberFormat = "@"
      $cell.Value = $entry[$c]
      $cell.Style = "Normal"
    }
  }
}

$wsPir = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPir $pirRows

$wsHum = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHum $humRows

